$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty inlineStr cells in column B of "ODI Batting" ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B4").ClearContents()
$batting.Range("B5").ClearContents()
$batting.Range("B8").ClearContents()

# --- 2. Add a new worksheet "ODI Batting Extra" at the end of the workbook ---
$extra = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$extra.Name = "ODI Batting Extra"

# Header row
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# MATCH_CODE values must stay text, so force text formatting before writing them
$extra.Range("A2:A10").NumberFormat = "@"

# Row 2 - match 4274
$extra.Range("A2").Value = "4274"
$extra.Range("B2").Value = 10
$extra.Range("F2").Value = "NO"

# Row 3 - match 4275
$extra.Range("A3").Value = "4275"
$extra.Range("F3").Value = "NO"

# Row 4 - match 4276
$extra.Range("A4").Value = "4276"
$extra.Range("B4").Value = 10
$extra.Range("F4").Value = "NO"

# Row 5 - match 4297
$extra.Range("A5").Value = "4297"
$extra.Range("B5").Value = 10
$extra.Range("F5").Value = "NO"

# Row 6 - match 4300
$extra.Range("A6").Value = "4300"
$extra.Range("F6").Value = "NO"

# Row 7 - match 4434
$extra.Range("A7").Value = "4434"

# Row 8 - match 4458
$extra.Range("A8").Value = "4458"

# Row 9 - match 4459
$extra.Range("A9").Value = "4459"

# Row 10 - match 4690
$extra.Range("A10").Value = "4690"

# Style the header row like the other sheets (bold, bordered, centered/top aligned)
$headerRange = $extra.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
